$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.465169666666667
$ws.Range("H2").Value = 4.395509000000001
$ws.Range("I2").Value = 0.03229814945245693
$ws.Range("J2").Value = 0.03229814945245692
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 103.4275383333333
$ws.Range("N2").Value = 310.282615
$ws.Range("O2").Value = 0.2485530285127421
$ws.Range("P2").Value = 0.2485530285127421
$ws.Range("Q2").Value = 151.5388918640039
$ws.Range("R2").Value = 1363.850026776035
$ws.Range("S2").Value = 0.008027802861765333
$ws.Range("T2").Value = 0.008027802861765332
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.465169666666667
$ws.Range("H3").Value = 4.395509000000001
$ws.Range("I3").Value = 0.03229814945245693
$ws.Range("J3").Value = 0.03229814945245692
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 216.130539
$ws.Range("N3").Value = 648.391617
$ws.Range("O3").Value = 0.5193964865470273
$ws.Range("P3").Value = 0.5193964865470272
$ws.Range("Q3").Value = 316.6679097831171
$ws.Range("R3").Value = 2850.011188048054
$ws.Range("S3").Value = 0.01677554534757692
$ws.Range("T3").Value = 0.01677554534757692
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.465169666666667
$ws.Range("H4").Value = 4.395509000000001
$ws.Range("I4").Value = 0.03229814945245693
$ws.Range("J4").Value = 0.03229814945245692
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 71.607325
$ws.Range("N4").Value = 214.821975
$ws.Range("O4").Value = 0.1720839321833696
$ws.Range("P4").Value = 0.1720839321833696
$ws.Range("Q4").Value = 104.9168805011417
$ws.Range("R4").Value = 944.2519245102751
$ws.Range("S4").Value = 0.005557992560024933
$ws.Range("T4").Value = 0.005557992560024932
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.465169666666667
$ws.Range("H5").Value = 4.395509000000001
$ws.Range("I5").Value = 0.03229814945245693
$ws.Range("J5").Value = 0.03229814945245692
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 24.953198
$ws.Range("N5").Value = 74.859594
$ws.Range("O5").Value = 0.05996655275686102
$ws.Range("P5").Value = 0.05996655275686102
$ws.Range("Q5").Value = 36.56066879592734
$ws.Range("R5").Value = 329.046019163346
$ws.Range("S5").Value = 0.001936808683089741
$ws.Range("T5").Value = 0.00193680868308974
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 40.16021866666667
$ws.Range("H6").Value = 120.480656
$ws.Range("I6").Value = 0.8852904711645572
$ws.Range("J6").Value = 0.8852904711645572
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 103.4275383333333
$ws.Range("N6").Value = 310.282615
$ws.Range("O6").Value = 0.2485530285127421
$ws.Range("P6").Value = 0.2485530285127421
$ws.Range("Q6").Value = 4153.672555621717
$ws.Range("R6").Value = 37383.05300059544
$ws.Range("S6").Value = 0.2200416277214231
$ws.Range("T6").Value = 0.220041627721423
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 40.16021866666667
$ws.Range("H7").Value = 120.480656
$ws.Range("I7").Value = 0.8852904711645572
$ws.Range("J7").Value = 0.8852904711645572
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 216.130539
$ws.Range("N7").Value = 648.391617
$ws.Range("O7").Value = 0.5193964865470273
$ws.Range("P7").Value = 0.5193964865470272
$ws.Range("Q7").Value = 8679.849706784529
$ws.Range("R7").Value = 78118.64736106076
$ws.Range("S7").Value = 0.4598167602964334
$ws.Range("T7").Value = 0.4598167602964333
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 40.16021866666667
$ws.Range("H8").Value = 120.480656
$ws.Range("I8").Value = 0.8852904711645572
$ws.Range("J8").Value = 0.8852904711645572
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 71.607325
$ws.Range("N8").Value = 214.821975
$ws.Range("O8").Value = 0.1720839321833696
$ws.Range("P8").Value = 0.1720839321833696
$ws.Range("Q8").Value = 2875.765830135067
$ws.Range("R8").Value = 25881.8924712156
$ws.Range("S8").Value = 0.152344265402465
$ws.Range("T8").Value = 0.152344265402465
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 40.16021866666667
$ws.Range("H9").Value = 120.480656
$ws.Range("I9").Value = 0.8852904711645572
$ws.Range("J9").Value = 0.8852904711645572
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 24.953198
$ws.Range("N9").Value = 74.859594
$ws.Range("O9").Value = 0.05996655275686102
$ws.Range("P9").Value = 0.05996655275686102
$ws.Range("Q9").Value = 1002.125888112629
$ws.Range("R9").Value = 9019.132993013665
$ws.Range("S9").Value = 0.05308781774423577
$ws.Range("T9").Value = 0.05308781774423577
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.279948333333333
$ws.Range("H10").Value = 3.839845
$ws.Range("I10").Value = 0.02821513678717742
$ws.Range("J10").Value = 0.02821513678717743
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 103.4275383333333
$ws.Range("N10").Value = 310.282615
$ws.Range("O10").Value = 0.2485530285127421
$ws.Range("P10").Value = 0.2485530285127421
$ws.Range("Q10").Value = 132.3819053105194
$ws.Range("R10").Value = 1191.437147794675
$ws.Range("S10").Value = 0.007012957698354229
$ws.Range("T10").Value = 0.007012957698354229
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1.279948333333333
$ws.Range("H11").Value = 3.839845
$ws.Range("I11").Value = 0.02821513678717742
$ws.Range("J11").Value = 0.02821513678717743
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 216.130539
$ws.Range("N11").Value = 648.391617
$ws.Range("O11").Value = 0.5193964865470273
$ws.Range("P11").Value = 0.5193964865470272
$ws.Range("Q11").Value = 276.635923175485
$ws.Range("R11").Value = 2489.723308579365
$ws.Range("S11").Value = 0.01465484291470373
$ws.Range("T11").Value = 0.01465484291470373
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 1.279948333333333
$ws.Range("H12").Value = 3.839845
$ws.Range("I12").Value = 0.02821513678717742
$ws.Range("J12").Value = 0.02821513678717743
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 71.607325
$ws.Range("N12").Value = 214.821975
$ws.Range("O12").Value = 0.1720839321833696
$ws.Range("P12").Value = 0.1720839321833696
$ws.Range("Q12").Value = 91.65367628820833
$ws.Range("R12").Value = 824.8830865938751
$ws.Range("S12").Value = 0.004855371685429135
$ws.Range("T12").Value = 0.004855371685429136
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 1.279948333333333
$ws.Range("H13").Value = 3.839845
$ws.Range("I13").Value = 0.02821513678717742
$ws.Range("J13").Value = 0.02821513678717743
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 24.953198
$ws.Range("N13").Value = 74.859594
$ws.Range("O13").Value = 0.05996655275686102
$ws.Range("P13").Value = 0.05996655275686102
$ws.Range("Q13").Value = 31.93880419143667
$ws.Range("R13").Value = 287.44923772293
$ws.Range("S13").Value = 0.001691964488690325
$ws.Range("T13").Value = 0.001691964488690325
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 2.458552333333333
$ws.Range("H14").Value = 7.375657
$ws.Range("I14").Value = 0.0541962425958086
$ws.Range("J14").Value = 0.0541962425958086
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 103.4275383333333
$ws.Range("N14").Value = 310.282615
$ws.Range("O14").Value = 0.2485530285127421
$ws.Range("P14").Value = 0.2485530285127421
$ws.Range("Q14").Value = 254.2820157003395
$ws.Range("R14").Value = 2288.538141303055
$ws.Range("S14").Value = 0.0134706402311995
$ws.Range("T14").Value = 0.0134706402311995
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 2.458552333333333
$ws.Range("H15").Value = 7.375657
$ws.Range("I15").Value = 0.0541962425958086
$ws.Range("J15").Value = 0.0541962425958086
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 216.130539
$ws.Range("N15").Value = 648.391617
$ws.Range("O15").Value = 0.5193964865470273
$ws.Range("P15").Value = 0.5193964865470272
$ws.Range("Q15").Value = 531.368240963041
$ws.Range("R15").Value = 4782.314168667369
$ws.Range("S15").Value = 0.02814933798831333
$ws.Range("T15").Value = 0.02814933798831332
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 2.458552333333333
$ws.Range("H16").Value = 7.375657
$ws.Range("I16").Value = 0.0541962425958086
$ws.Range("J16").Value = 0.0541962425958086
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 71.607325
$ws.Range("N16").Value = 214.821975
$ws.Range("O16").Value = 0.1720839321833696
$ws.Range("P16").Value = 0.1720839321833696
$ws.Range("Q16").Value = 176.0503559625083
$ws.Range("R16").Value = 1584.453203662575
$ws.Range("S16").Value = 0.009326302535450573
$ws.Range("T16").Value = 0.009326302535450573
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 2.458552333333333
$ws.Range("H17").Value = 7.375657
$ws.Range("I17").Value = 0.0541962425958086
$ws.Range("J17").Value = 0.0541962425958086
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 24.953198
$ws.Range("N17").Value = 74.859594
$ws.Range("O17").Value = 0.05996655275686102
$ws.Range("P17").Value = 0.05996655275686102
$ws.Range("Q17").Value = 61.34874316702867
$ws.Range("R17").Value = 552.138688503258
$ws.Range("S17").Value = 0.003249961840845195
$ws.Range("T17").Value = 0.003249961840845195
